# Update cryptocurrency price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '27.045.92'
$ws.Range('E2').Value = '  +0.68%  '

$ws.Range('D3').Value = '1.683.49'
$ws.Range('E3').Value = '  +0.89%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.28%  '

$ws.Range('E6').Value = '  -2.27%  '

$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.59'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.72%  '

$ws.Range('E9').Value = '  -0.18%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0622'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.57%  '

$ws.Range('E11').Value = '  -0.43%  '

$ws.Range('D12').Value = '1.921.23'
$ws.Range('E12').Value = '  +0.89%  '

$ws.Range('D13').Value = '1.680.94'
$ws.Range('E13').Value = '  +0.39%  '

$ws.Range('E14').Value = '  +0.50%  '

$ws.Range('E15').Value = '  +1.89%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.86%  '

$ws.Range('D17').Value = '27.050.51'
$ws.Range('E17').Value = '  +0.58%  '

$ws.Range('E18').Value = '  +5.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '236.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.12%  '

$ws.Range('E20').Value = '  +0.51%  '

$ws.Range('E21').Value = '  -0.02%  '

$ws.Range('E22').Value = '  +0.14%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.78%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.99%  '

$ws.Range('E26').Value = '  +5.68%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.78%  '

$ws.Range('E28').Value = '  -2.78%  '

$ws.Range('E29').Value = '  +0.18%  '

$ws.Range('E30').Value = '  +0.48%  '

$ws.Range('E31').Value = '  -0.15%  '

$ws.Range('E32').Value = '  +0.37%  '

$ws.Range('D33').Value = '1.525.16'
$ws.Range('E33').Value = '  +4.21%  '

$ws.Range('E34').Value = '  +1.03%  '

$ws.Range('E35').Value = '  +4.64%  '

$ws.Range('E36').Value = '  -0.56%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.591'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.64%  '

$ws.Range('E38').Value = '  +2.63%  '

$ws.Range('E39').Value = '  +3.75%  '

$ws.Range('E40').Value = '  +6.27%  '

$ws.Range('E41').Value = '  -0.89%  '

$ws.Range('E42').Value = '  -0.01%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '68.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.67%  '

$ws.Range('D45').Value = '1.825.50'
$ws.Range('E45').Value = '  +0.54%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.782'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.12%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.28'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.02%  '

$ws.Range('E48').Value = '  -0.26%  '

$ws.Range('E49').Value = '  +4.23%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.18%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.43%  '
